$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Helper: write a date-like label into a cell as plain text (shared
# string) instead of letting Excel auto-convert it to a date serial.
# We enter it with a leading apostrophe (text-qualifier) and then strip
# the resulting "quote prefix" cell formatting back off so the cell ends
# up using the workbook's default (unstyled) cell format, matching the
# rest of the sheet.
# ----------------------------------------------------------------------
function Set-TextValue($ws, $row, $col, $text) {
    $ws.Cells.Item($row, $col).Value = "'" + $text
}

$newDates = @("04/02/2022","05/02/2022","06/02/2022","07/02/2022","08/02/2022","09/02/2022","10/02/2022")

# ======================================================================
# Sheet "ASKARI" (sheet1)
# ======================================================================
$ws1 = $wb.Worksheets.Item("ASKARI")

# Rows 53-65: loan balance increased -> H changes from 30,000,000 to
# 130,000,000 and the markup (J) changes accordingly.
for ($r = 53; $r -le 65; $r++) {
    $ws1.Cells.Item($r, 8).Value = 130000000
    $ws1.Cells.Item($r, 10).Value = 80087.67123287672
}

# Row 53 also gains a new disbursement recorded in E53/F53.
$ws1.Cells.Item(53, 5).Value = 12.9
$ws1.Cells.Item(53, 6).Value = 100000000

# Rows 66-74: loan balance increased again -> H changes from 30,000,000
# to 230,000,000 and the markup (J) changes accordingly.
for ($r = 66; $r -le 74; $r++) {
    $ws1.Cells.Item($r, 8).Value = 230000000
    $ws1.Cells.Item($r, 10).Value = 80087.67123287672
}

# Row 66 gains a new disbursement recorded in E66/F66.
$ws1.Cells.Item(66, 5).Value = 12.55
$ws1.Cells.Item(66, 6).Value = 100000000

# The disbursement previously recorded on row 72 (E72/F72) no longer
# applies -- it moved to row 53/66 above -- so clear it out.
$ws1.Range("E72:F72").ClearContents()

# New rows 75-81 continuing the amortisation schedule.
$row75_81 = @(
    @{ B = 73; L = 9402.739726027396 },
    @{ B = 74; L = 18805.47945205479 },
    @{ B = 75; L = 28208.21917808219 },
    @{ B = 76; L = 37610.95890410958 },
    @{ B = 77; L = 47013.69863013698 },
    @{ B = 78; L = 56416.43835616437 },
    @{ B = 79; L = 80087.67123287672 }
)
for ($i = 0; $i -lt 7; $i++) {
    $r = 75 + $i
    $ws1.Cells.Item($r, 2).Value = $row75_81[$i].B
    Set-TextValue $ws1 $r 3 $newDates[$i]
    $ws1.Cells.Item($r, 8).Value = 230000000
    $ws1.Cells.Item($r, 10).Value = 80087.67123287672
    $ws1.Cells.Item($r, 12).Value = $row75_81[$i].L
}
$ws1.Range("C75:C81").ClearFormats()

# ======================================================================
# Sheet "DIBL" (sheet2)
# ======================================================================
$ws2 = $wb.Worksheets.Item("DIBL")

$row41_47 = @(
    @{ B = 39; L = 16397.2602739726 },
    @{ B = 40; L = 32794.52054794521 },
    @{ B = 41; L = 49191.78082191781 },
    @{ B = 42; L = 65589.04109589041 },
    @{ B = 43; L = 81986.30136986301 },
    @{ B = 44; L = 98383.56164383561 },
    @{ B = 45; L = 16397.2602739726 }
)
for ($i = 0; $i -lt 7; $i++) {
    $r = 41 + $i
    $ws2.Cells.Item($r, 2).Value = $row41_47[$i].B
    Set-TextValue $ws2 $r 3 $newDates[$i]
    $ws2.Cells.Item($r, 8).Value = 45000000
    $ws2.Cells.Item($r, 10).Value = 16397.2602739726
    $ws2.Cells.Item($r, 12).Value = $row41_47[$i].L
}
$ws2.Range("C41:C47").ClearFormats()

# ======================================================================
# Sheet "HBL" (sheet3)
# ======================================================================
$ws3 = $wb.Worksheets.Item("HBL")

$row69_75 = @(
    @{ B = 67; L = 9285.205479452054 },
    @{ B = 68; L = 18570.41095890411 },
    @{ B = 69; L = 27855.61643835616 },
    @{ B = 70; L = 37140.82191780821 },
    @{ B = 71; L = 46426.02739726027 },
    @{ B = 72; L = 55711.23287671232 },
    @{ B = 73; L = 9285.205479452054 }
)
for ($i = 0; $i -lt 7; $i++) {
    $r = 69 + $i
    $ws3.Cells.Item($r, 2).Value = $row69_75[$i].B
    Set-TextValue $ws3 $r 3 $newDates[$i]
    $ws3.Cells.Item($r, 8).Value = 28600000
    $ws3.Cells.Item($r, 10).Value = 9285.205479452054
    $ws3.Cells.Item($r, 12).Value = $row69_75[$i].L
}
$ws3.Range("C69:C75").ClearFormats()
